# Removing less than USD 5 price from extrapolation calibration because it is just a noise
# This updates the recalculated calibration outputs (ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN)
# for the rows affected by dropping the sub-$5 option price from the fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 124178.1998038432
$ws.Range("E3").Value = 0.02748126053704688
$ws.Range("F3").Value = 0.1834936993069084
$ws.Range("G3").Value = -0.6204895358439093
$ws.Range("H3").Value = 9.471443401656614

$ws.Range("D5").Value = 125918.4996449818
$ws.Range("E5").Value = 0.005433370494172475
$ws.Range("F5").Value = 0.2194350031705524
$ws.Range("G5").Value = -0.7010986018052038
$ws.Range("H5").Value = 8.383799501217444

$ws.Range("D6").Value = 126359.6383706355
$ws.Range("E6").Value = -0.009663359226640628
$ws.Range("F6").Value = 0.2555379863949271
$ws.Range("G6").Value = -1.222310607217036
$ws.Range("H6").Value = 11.11734655883518

$ws.Range("D9").Value = 129636.1709660843
$ws.Range("E9").Value = -0.07040383814255735
$ws.Range("F9").Value = 0.3657060939354542
$ws.Range("G9").Value = -1.851251768435398
$ws.Range("H9").Value = 11.87399375507887

$ws.Range("D10").Value = 131128.8835390481
$ws.Range("E10").Value = -0.0940073360070954
$ws.Range("F10").Value = 0.4008313968713341
$ws.Range("G10").Value = -1.802891082623549
$ws.Range("H10").Value = 9.838719582597273

$ws.Range("D11").Value = 133308.7910640123
$ws.Range("E11").Value = -0.1681830936976825
$ws.Range("F11").Value = 0.7175581915692382
$ws.Range("G11").Value = -2.643298237415924
$ws.Range("H11").Value = 13.19983396292374

$ws.Range("D13").Value = 123419.2413078245
$ws.Range("E13").Value = 0.1004679077765962
$ws.Range("F13").Value = 0.1615421857277089
$ws.Range("G13").Value = -0.2592978026256869
$ws.Range("H13").Value = 7.053690629184262

$ws.Range("D14").Value = 123462.6316935806
$ws.Range("E14").Value = 0.09346904324945446
$ws.Range("F14").Value = 0.1672319128154472
$ws.Range("G14").Value = -0.06909789249204973
$ws.Range("H14").Value = 7.262250765205483
